$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B7: was inline string " 100", now numeric 80
$ws.Range("B7").Value = 80

# New row 8: A8 = "Teste" (merged A8:B8)
$ws.Range("A8").Value = "Teste"
$ws.Range("A8:B8").Merge()

# New row 9: B9 = 1280 (sum of expenses)
$ws.Range("B9").Value = 1280
